$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")

$row = 31
$ws.Cells.Item($row, 1).Value = "General Knowledge"
$ws.Cells.Item($row, 2).Value = "Sanjib Roy"
$ws.Cells.Item($row, 3).Value = "sanjibrosnjssjnjy0098@gmail.com"
$ws.Cells.Item($row, 4).Value = 2
$ws.Cells.Item($row, 5).Value = 0
